# edit.ps1 - reproduce the target commit:
#  1) Slide 16's table (graphicFrame "Google Shape;213;p29") switches its
#     table style from {3E89479C-4CD6-4B38-9310-7E9A4C750CB9} to
#     {6F4C8049-5532-46CF-BB2F-2AF9EE2C3138}.
#  2) The deck's theme (the one applied to the slide master, i.e. the
#     "Design") is swapped from the custom "Integral" colour set to the
#     stock "Office Theme" colour set (12 theme colours: dk1, lt1, dk2,
#     lt2, accent1-6, hlink, folHlink).

$p = $ppt.ActivePresentation

# --- 1) Table style on slide 16 -------------------------------------------
$s16 = $p.Slides.Item(16)
$tableShape = $s16.Shapes.Item(3)
$tableShape.Table.ApplyStyle("{6F4C8049-5532-46CF-BB2F-2AF9EE2C3138}")

# --- 2) Swap the theme colour scheme back to the stock Office colours ----
# (the slide master's theme part currently holds the "Integral" palette;
#  repoint every theme colour slot at the classic Office Theme palette)
$tcs = $s16.ThemeColorScheme
$tcs.Item(1).RGB  = 0          # dk1      -> 000000
$tcs.Item(2).RGB  = 16777215   # lt1      -> FFFFFF
$tcs.Item(3).RGB  = 6968388    # dk2      -> 44546A
$tcs.Item(4).RGB  = 15132391   # lt2      -> E7E6E6
$tcs.Item(5).RGB  = 13998939   # accent1  -> 5B9BD5
$tcs.Item(6).RGB  = 3243501    # accent2  -> ED7D31
$tcs.Item(7).RGB  = 10855845   # accent3  -> A5A5A5
$tcs.Item(8).RGB  = 49407      # accent4  -> FFC000
$tcs.Item(9).RGB  = 12874308   # accent5  -> 4472C4
$tcs.Item(10).RGB = 4697456    # accent6  -> 70AD47
$tcs.Item(11).RGB = 12673797   # hlink    -> 0563C1
$tcs.Item(12).RGB = 7491477    # folHlink -> 954F72
